# Apply updated NATMI edge-weight computations (Lrpap1-Sort1) following Dr Hou advice.
# Ligand/Receptor-expressing cell counts increase from 1 to 3 per cluster,
# which changes total/average expression values and derived specificities.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 6.660188000000001
$ws.Range("H2").Value = 19.980564
$ws.Range("I2").Value = 0.1500148400131262
$ws.Range("J2").Value = 0.1500148400131261
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 3.024839333333334
$ws.Range("N2").Value = 9.074518000000001
$ws.Range("O2").Value = 0.1801507982970389
$ws.Range("P2").Value = 0.1801507982970388
$ws.Range("Q2").Value = 20.14599862979467
$ws.Range("R2").Value = 181.313987668152
$ws.Range("S2").Value = 0.02702529318476725
$ws.Range("T2").Value = 0.02702529318476724

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 6.660188000000001
$ws.Range("H3").Value = 19.980564
$ws.Range("I3").Value = 0.1500148400131262
$ws.Range("J3").Value = 0.1500148400131261
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 1.626140333333333
$ws.Range("N3").Value = 4.878420999999999
$ws.Range("O3").Value = 0.09684827751501936
$ws.Range("P3").Value = 0.09684827751501934
$ws.Range("Q3").Value = 10.83040033438267
$ws.Range("R3").Value = 97.47360300944399
$ws.Range("S3").Value = 0.01452867885696247
$ws.Range("T3").Value = 0.01452867885696247

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 6.660188000000001
$ws.Range("H4").Value = 19.980564
$ws.Range("I4").Value = 0.1500148400131262
$ws.Range("J4").Value = 0.1500148400131261
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 4.640628666666666
$ws.Range("N4").Value = 13.921886
$ws.Range("O4").Value = 0.2763825997921178
$ws.Range("P4").Value = 0.2763825997921177
$ws.Range("Q4").Value = 30.90745935818934
$ws.Range("R4").Value = 278.167134223704
$ws.Range("S4").Value = 0.04146149149022643
$ws.Range("T4").Value = 0.04146149149022641

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 6.660188000000001
$ws.Range("H5").Value = 19.980564
$ws.Range("I5").Value = 0.1500148400131262
$ws.Range("J5").Value = 0.1500148400131261
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 7.498988000000001
$ws.Range("N5").Value = 22.496964
$ws.Range("O5").Value = 0.4466183243958241
$ws.Range("P5").Value = 0.446618324395824
$ws.Range("Q5").Value = 49.94466988974401
$ws.Range("R5").Value = 449.5020290076961
$ws.Range("S5").Value = 0.06699937648117003
$ws.Range("T5").Value = 0.06699937648117

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 16.410331
$ws.Range("H6").Value = 49.230993
$ws.Range("I6").Value = 0.3696281815959916
$ws.Range("J6").Value = 0.3696281815959916
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 3.024839333333334
$ws.Range("N6").Value = 9.074518000000001
$ws.Range("O6").Value = 0.1801507982970389
$ws.Range("P6").Value = 0.1801507982970388
$ws.Range("Q6").Value = 49.63861468181934
$ws.Range("R6").Value = 446.747532136374
$ws.Range("S6").Value = 0.06658881198760075
$ws.Range("T6").Value = 0.06658881198760074

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 16.410331
$ws.Range("H7").Value = 49.230993
$ws.Range("I7").Value = 0.3696281815959916
$ws.Range("J7").Value = 0.3696281815959916
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 1.626140333333333
$ws.Range("N7").Value = 4.878420999999999
$ws.Range("O7").Value = 0.09684827751501936
$ws.Range("P7").Value = 0.09684827751501934
$ws.Range("Q7").Value = 26.68550112245033
$ws.Range("R7").Value = 240.169510102053
$ws.Range("S7").Value = 0.03579785270858057
$ws.Range("T7").Value = 0.03579785270858057

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 16.410331
$ws.Range("H8").Value = 49.230993
$ws.Range("I8").Value = 0.3696281815959916
$ws.Range("J8").Value = 0.3696281815959916
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 4.640628666666666
$ws.Range("N8").Value = 13.921886
$ws.Range("O8").Value = 0.2763825997921178
$ws.Range("P8").Value = 0.2763825997921177
$ws.Range("Q8").Value = 76.15425246808866
$ws.Range("R8").Value = 685.3882722127979
$ws.Range("S8").Value = 0.1021587977859332
$ws.Range("T8").Value = 0.1021587977859332

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 16.410331
$ws.Range("H9").Value = 49.230993
$ws.Range("I9").Value = 0.3696281815959916
$ws.Range("J9").Value = 0.3696281815959916
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 7.498988000000001
$ws.Range("N9").Value = 22.496964
$ws.Range("O9").Value = 0.4466183243958241
$ws.Range("P9").Value = 0.446618324395824
$ws.Range("Q9").Value = 123.060875245028
$ws.Range("R9").Value = 1107.547877205252
$ws.Range("S9").Value = 0.1650827191138771
$ws.Range("T9").Value = 0.1650827191138771

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 14.478895
$ws.Range("H10").Value = 43.436685
$ws.Range("I10").Value = 0.3261242951387937
$ws.Range("J10").Value = 0.3261242951387937
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 3.024839333333334
$ws.Range("N10").Value = 9.074518000000001
$ws.Range("O10").Value = 0.1801507982970389
$ws.Range("P10").Value = 0.1801507982970388
$ws.Range("Q10").Value = 43.79633109920334
$ws.Range("R10").Value = 394.16697989283
$ws.Range("S10").Value = 0.0587515521133128
$ws.Range("T10").Value = 0.05875155211331278

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 14.478895
$ws.Range("H11").Value = 43.436685
$ws.Range("I11").Value = 0.3261242951387937
$ws.Range("J11").Value = 0.3261242951387937
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 1.626140333333333
$ws.Range("N11").Value = 4.878420999999999
$ws.Range("O11").Value = 0.09684827751501936
$ws.Range("P11").Value = 0.09684827751501934
$ws.Range("Q11").Value = 23.54471514159833
$ws.Range("R11").Value = 211.902436274385
$ws.Range("S11").Value = 0.03158457623999197
$ws.Range("T11").Value = 0.03158457623999196

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 14.478895
$ws.Range("H12").Value = 43.436685
$ws.Range("I12").Value = 0.3261242951387937
$ws.Range("J12").Value = 0.3261242951387937
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 4.640628666666666
$ws.Range("N12").Value = 13.921886
$ws.Range("O12").Value = 0.2763825997921178
$ws.Range("P12").Value = 0.2763825997921177
$ws.Range("Q12").Value = 67.19117519865667
$ws.Range("R12").Value = 604.7205767879099
$ws.Range("S12").Value = 0.09013508054583172
$ws.Range("T12").Value = 0.09013508054583169

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 14.478895
$ws.Range("H13").Value = 43.436685
$ws.Range("I13").Value = 0.3261242951387937
$ws.Range("J13").Value = 0.3261242951387937
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 7.498988000000001
$ws.Range("N13").Value = 22.496964
$ws.Range("O13").Value = 0.4466183243958241
$ws.Range("P13").Value = 0.446618324395824
$ws.Range("Q13").Value = 108.57705985826
$ws.Range("R13").Value = 977.19353872434
$ws.Range("S13").Value = 0.1456530862396572
$ws.Range("T13").Value = 0.1456530862396572

# Row 14
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 6.847447
$ws.Range("H14").Value = 20.542341
$ws.Range("I14").Value = 0.1542326832520885
$ws.Range("J14").Value = 0.1542326832520885
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 3.024839333333334
$ws.Range("N14").Value = 9.074518000000001
$ws.Range("O14").Value = 0.1801507982970389
$ws.Range("P14").Value = 0.1801507982970388
$ws.Range("Q14").Value = 20.71242701851533
$ws.Range("R14").Value = 186.411843166638
$ws.Range("S14").Value = 0.02778514101135807
$ws.Range("T14").Value = 0.02778514101135807

# Row 15
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 6.847447
$ws.Range("H15").Value = 20.542341
$ws.Range("I15").Value = 0.1542326832520885
$ws.Range("J15").Value = 0.1542326832520885
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 1.626140333333333
$ws.Range("N15").Value = 4.878420999999999
$ws.Range("O15").Value = 0.09684827751501936
$ws.Range("P15").Value = 0.09684827751501934
$ws.Range("Q15").Value = 11.13490974706233
$ws.Range("R15").Value = 100.214187723561
$ws.Range("S15").Value = 0.01493716970948434
$ws.Range("T15").Value = 0.01493716970948434

# Row 16
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 6.847447
$ws.Range("H16").Value = 20.542341
$ws.Range("I16").Value = 0.1542326832520885
$ws.Range("J16").Value = 0.1542326832520885
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 4.640628666666666
$ws.Range("N16").Value = 13.921886
$ws.Range("O16").Value = 0.2763825997921178
$ws.Range("P16").Value = 0.2763825997921177
$ws.Range("Q16").Value = 31.77645884168066
$ws.Range("R16").Value = 285.988129575126
$ws.Range("S16").Value = 0.04262722997012643
$ws.Range("T16").Value = 0.04262722997012643

# Row 17
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 6.847447
$ws.Range("H17").Value = 20.542341
$ws.Range("I17").Value = 0.1542326832520885
$ws.Range("J17").Value = 0.1542326832520885
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 7.498988000000001
$ws.Range("N17").Value = 22.496964
$ws.Range("O17").Value = 0.4466183243958241
$ws.Range("P17").Value = 0.446618324395824
$ws.Range("Q17").Value = 51.348922883636
$ws.Range("R17").Value = 462.14030595272396
$ws.Range("S17").Value = 0.06888314256111965
$ws.Range("T17").Value = 0.06888314256111963

